# Refresh the cryptos price/volume table (GitHub Actions data pull).
# For D-column prices that parse as a plain decimal number, an apostrophe
# prefix forces Excel to keep them as text (matching the sheet's original
# inlineStr cells) and .Style = "Normal" strips the resulting quote-prefix
# formatting back off so no stray number format is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.761.33"
$ws.Range("E2").Value = "  +4.09%  "

$ws.Range("D3").Value = "2.266.49"
$ws.Range("E3").Value = "  +2.08%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'304.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.09%  "

$ws.Range("D6").Value = "'91.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.14%  "

$ws.Range("D7").Value = "'0.530"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.26%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  +1.92%  "

$ws.Range("E10").Value = "  +4.78%  "

$ws.Range("D11").Value = "'53.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.12%  "

$ws.Range("D12").Value = "'0.0794"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.47%  "

$ws.Range("E13").Value = "  +0.82%  "

$ws.Range("D14").Value = "'6.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.66%  "

$ws.Range("D15").Value = "2.617.48"
$ws.Range("E15").Value = "  +2.13%  "

$ws.Range("E16").Value = "  +1.91%  "

$ws.Range("D17").Value = "2.290.46"
$ws.Range("E17").Value = "  +2.22%  "

$ws.Range("E18").Value = "  +3.12%  "

$ws.Range("D19").Value = "41.705.46"
$ws.Range("E19").Value = "  +4.11%  "

$ws.Range("D20").Value = "'12.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.18%  "

$ws.Range("E21").Value = "  +1.55%  "

$ws.Range("D22").Value = "'5.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.03%  "

$ws.Range("D23").Value = "'66.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.52%  "

$ws.Range("D24").Value = "'241.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.53%  "

$ws.Range("E25").Value = "  +3.69%  "

$ws.Range("E27").Value = "  +4.34%  "

$ws.Range("D28").Value = "'24.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.10%  "

$ws.Range("D29").Value = "'2.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +10.41%  "

$ws.Range("D30").Value = "'9.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.79%  "

$ws.Range("D31").Value = "'159.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.01%  "

$ws.Range("D32").Value = "'34.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.70%  "

$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("E34").Value = "  +3.74%  "

$ws.Range("D35").Value = "'0.0742"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.90%  "

$ws.Range("D36").Value = "'3.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.03%  "

$ws.Range("E37").Value = "  +2.01%  "

$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").Value = "'0.116"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.47%  "

$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").Value = "'16.61"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.16%  "

$ws.Range("E40").Value = "  +3.34%  "

$ws.Range("E41").Value = "  +1.79%  "

$ws.Range("E42").Value = "  +3.92%  "

$ws.Range("D43").Value = "2.059.14"
$ws.Range("E43").Value = "  -1.05%  "

$ws.Range("D44").Value = "'19.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.20%  "

$ws.Range("E45").Value = "  +2.41%  "

$ws.Range("D46").Value = "'10.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.11%  "

$ws.Range("E47").Value = "  +3.00%  "

$ws.Range("E48").Value = "  +6.48%  "

$ws.Range("D49").Value = "'73.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.02%  "

$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").Value = "'1.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.33%  "

$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'1.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.25%  "
